$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 375.92856
$ws.Range("I19").Value = 219.75
$ws.Range("J19").Value = 438.4
$ws.Range("K19").Value = 219.75
$ws.Range("L19").Value = 438.4
$ws.Range("M19").Value = -44.75
$ws.Range("N19").Value = -788.4

$ws.Range("H113").Value = 9526409
$ws.Range("I113").Value = 11113661
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 11113661
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -11110407
$ws.Range("N113").Value = -9408

$ws.Range("H138").Value = 1721.2424
$ws.Range("I138").Value = 836.3871
$ws.Range("J138").Value = 2504.9714
$ws.Range("K138").Value = 2509.1613
$ws.Range("L138").Value = 7514.914199999999
$ws.Range("M138").Value = 2630.8387
$ws.Range("N138").Value = -17794.9142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6628.846
$ws.Range("I32").Value = 4663.716
$ws.Range("K32").Value = 4663.716
$ws.Range("M32").Value = -4376.716

$ws.Range("H61").Value = 215440.89
$ws.Range("I61").Value = 1921.5927
$ws.Range("K61").Value = 1921.5927
$ws.Range("M61").Value = -1709.5927

$ws.Range("H132").Value = 5004577
$ws.Range("I132").Value = 2902.9092
$ws.Range("J132").Value = 11117735
$ws.Range("K132").Value = 8708.7276
$ws.Range("L132").Value = 33353205
$ws.Range("M132").Value = -6178.7276
$ws.Range("N132").Value = -33358265

$ws.Range("H136").Value = 215440.89
$ws.Range("I136").Value = 1921.5927
$ws.Range("K136").Value = 5764.7781
$ws.Range("M136").Value = -3214.7781

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 24828.818
$ws.Range("I20").Value = 2143.4
$ws.Range("J20").Value = 43733.332
$ws.Range("K20").Value = 2143.4
$ws.Range("L20").Value = 43733.332
$ws.Range("M20").Value = -1896.4
$ws.Range("N20").Value = -44227.332

$ws.Range("H134").Value = 2600.6296
$ws.Range("I134").Value = 2183.9092
$ws.Range("K134").Value = 6551.7276
$ws.Range("M134").Value = -4016.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 571.8
$ws.Range("I22").Value = 667.25
$ws.Range("J22").Value = 190
$ws.Range("K22").Value = 667.25
$ws.Range("L22").Value = 190
$ws.Range("M22").Value = -317.25
$ws.Range("N22").Value = -890

$ws.Range("H44").Value = 500004450
$ws.Range("J44").Value = 500004450
$ws.Range("L44").Value = 500004450
$ws.Range("N44").Value = -500005334

$ws.Range("H94").Value = 3554.8057
$ws.Range("I94").Value = 2649.4119
$ws.Range("J94").Value = 4364.8945
$ws.Range("K94").Value = 2649.4119
$ws.Range("L94").Value = 4364.8945
$ws.Range("M94").Value = -2198.4119
$ws.Range("N94").Value = -5266.8945

$ws.Range("H99").Value = 6629.8
$ws.Range("I99").Value = 7924.75
$ws.Range("J99").Value = 1450
$ws.Range("K99").Value = 7924.75
$ws.Range("L99").Value = 1450
$ws.Range("M99").Value = -6426.75
$ws.Range("N99").Value = -4446

$ws.Range("H122").Value = 2318151
$ws.Range("I122").Value = 5556351
$ws.Range("K122").Value = 16669053
$ws.Range("M122").Value = -16666603

$ws.Range("H126").Value = 6629.8
$ws.Range("I126").Value = 7924.75
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 23774.25
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -21304.25
$ws.Range("N126").Value = -9290

$ws.Range("H134").Value = 2468.1667
$ws.Range("I134").Value = 2469.4736
$ws.Range("J134").Value = 2455.75
$ws.Range("K134").Value = 7408.4208
$ws.Range("L134").Value = 7367.25
$ws.Range("M134").Value = -4873.4208
$ws.Range("N134").Value = -12437.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 2123.4443
$ws.Range("I46").Value = 1111
$ws.Range("J46").Value = 2250
$ws.Range("K46").Value = 3333
$ws.Range("L46").Value = 6750
$ws.Range("M46").Value = -3242
$ws.Range("N46").Value = -6932

$ws.Range("H112").Value = 2909.5454
$ws.Range("J112").Value = 3046.0977
$ws.Range("L112").Value = 9138.293099999999
$ws.Range("N112").Value = -11354.2931

$ws.Range("H132").Value = 1491.7273
$ws.Range("I132").Value = 584
$ws.Range("J132").Value = 2581
$ws.Range("K132").Value = 5256
$ws.Range("L132").Value = 23229
$ws.Range("M132").Value = -2726
$ws.Range("N132").Value = -28289

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4964.2104
$ws.Range("I70").Value = 4984.6787
$ws.Range("J70").Value = 4906.9
$ws.Range("K70").Value = 4984.6787
$ws.Range("L70").Value = 4906.9
$ws.Range("M70").Value = -4714.6787
$ws.Range("N70").Value = -5446.9

$ws.Range("H73").Value = 4964.2104
$ws.Range("I73").Value = 4984.6787
$ws.Range("J73").Value = 4906.9
$ws.Range("K73").Value = 4984.6787
$ws.Range("L73").Value = 4906.9
$ws.Range("M73").Value = -4048.6787
$ws.Range("N73").Value = -6778.9

$ws.Range("H132").Value = 4099.2856
$ws.Range("I132").Value = 5552.5557
$ws.Range("J132").Value = 3009.3333
$ws.Range("K132").Value = 16657.6671
$ws.Range("L132").Value = 9027.999899999999
$ws.Range("M132").Value = -14127.6671
$ws.Range("N132").Value = -14087.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 100101864
$ws.Range("I68").Value = 144664.72
$ws.Range("J68").Value = 333335330
$ws.Range("K68").Value = 144664.72
$ws.Range("L68").Value = 333335330
$ws.Range("M68").Value = -143915.72
$ws.Range("N68").Value = -333336828

$ws.Range("H71").Value = 100101864
$ws.Range("I71").Value = 144664.72
$ws.Range("J71").Value = 333335330
$ws.Range("K71").Value = 723323.6
$ws.Range("L71").Value = 1666676650
$ws.Range("M71").Value = -719579.6
$ws.Range("N71").Value = -1666684138

$ws.Range("H132").Value = 18522230
$ws.Range("I132").Value = 27780604
$ws.Range("J132").Value = 5483
$ws.Range("K132").Value = 83341812
$ws.Range("L132").Value = 16449
$ws.Range("M132").Value = -83339282
$ws.Range("N132").Value = -21509

$ws.Range("H136").Value = 3682.1292
$ws.Range("I136").Value = 1964.375
$ws.Range("J136").Value = 9571.571
$ws.Range("K136").Value = 5893.125
$ws.Range("L136").Value = 28714.713
$ws.Range("M136").Value = -3343.125
$ws.Range("N136").Value = -33814.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3847.3333
$ws.Range("I62").Value = 3847.3333
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3847.3333
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3223.3333
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3847.3333
$ws.Range("I65").Value = 3847.3333
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 19236.6665
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16116.6665
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 861.93335
$ws.Range("I122").Value = 711.2727
$ws.Range("J122").Value = 1276.25
$ws.Range("K122").Value = 2133.8181
$ws.Range("L122").Value = 3828.75
$ws.Range("M122").Value = 316.1819
$ws.Range("N122").Value = -8728.75

$ws.Range("H126").Value = 772.3125
$ws.Range("I126").Value = 596.9286
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 1790.7858
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = 679.2142000000001
$ws.Range("N126").Value = -10940

$ws.Range("H132").Value = 1574.5227
$ws.Range("I132").Value = 1151.2142
$ws.Range("J132").Value = 2315.3125
$ws.Range("K132").Value = 3453.6426
$ws.Range("L132").Value = 6945.9375
$ws.Range("M132").Value = -923.6425999999997
$ws.Range("N132").Value = -12005.9375

$ws.Range("H136").Value = 2662.946
$ws.Range("I136").Value = 3433.375
$ws.Range("J136").Value = 2075.9524
$ws.Range("K136").Value = 10300.125
$ws.Range("L136").Value = 6227.8572
$ws.Range("M136").Value = -7750.125
$ws.Range("N136").Value = -11327.8572
